$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 77267

# Row 3
$ws.Range("B3").Value = 96348

# Row 4
$ws.Range("A4").Value = 111052330
$ws.Range("B4").Value = 96348
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "21"
$ws.Range("J4").Value = "plantor/tuvor"
$ws.Range("P4").Value = "Ön Haverö (Ön Haverö), Mpd"
$ws.Range("Q4").Value = 505566.4405733302
$ws.Range("R4").Value = 6913047.682812326
$ws.Range("S4").Value = 20
$ws.Range("Z4").Value = "11:52"
$ws.Range("AB4").Value = "11:52"

# Row 5
$ws.Range("A5").Value = 111052279
$ws.Range("B5").Value = 96348
$ws.Range("D5").Value = "VU"
$ws.Range("E5").Value = 220787
$ws.Range("F5").Value = "Knärot"
$ws.Range("G5").Value = "Goodyera repens"
$ws.Range("H5").Value = "(L.) R. Br."
$ws.Range("Q5").Value = 505548.3107360627
$ws.Range("R5").Value = 6913041.614857432
$ws.Range("S5").Value = 20
$ws.Range("Z5").Value = "11:44"
$ws.Range("AB5").Value = "11:44"

# Row 6
$ws.Range("A6").Value = 111052223
$ws.Range("B6").Value = 96348
$ws.Range("P6").Value = "Ömansholmen (Ömansholmen), Mpd"
$ws.Range("Q6").Value = 505585.5400955554
$ws.Range("R6").Value = 6913030.065418502

# Row 7
$ws.Range("A7").Value = 111052010
$ws.Range("B7").Value = 77268
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 228912
$ws.Range("F7").Value = "Mörk kolflarnlav"
$ws.Range("G7").Value = "Carbonicola myrmecina"
$ws.Range("H7").Value = "(Ach.) Bendiksby & Timdal"
$ws.Range("I7").Value = ""
$ws.Range("J7").Value = ""
$ws.Range("P7").Value = "Ömansholmen (Ömansholmen), Mpd"
$ws.Range("Q7").Value = 505537.2690634067
$ws.Range("R7").Value = 6912968.212359355
$ws.Range("Z7").Value = "11:15"
$ws.Range("AB7").Value = "11:15"

# Row 8
$ws.Range("A8").Value = 111052497
$ws.Range("B8").Value = 77267
$ws.Range("E8").Value = 6446
$ws.Range("F8").Value = "Kolflarnlav"
$ws.Range("G8").Value = "Carbonicola anthracophila"
$ws.Range("H8").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("P8").Value = "Ön Haverö (Ön Haverö), Mpd"
$ws.Range("Q8").Value = 505544.8695405752
$ws.Range("R8").Value = 6913153.540235949
$ws.Range("S8").Value = 10
$ws.Range("Z8").Value = "11:52"
$ws.Range("AB8").Value = "11:52"

# Row 9
$ws.Range("A9").Value = 111052144
$ws.Range("B9").Value = 96348
$ws.Range("Q9").Value = 505586.595854046
$ws.Range("R9").Value = 6912954.824843475
$ws.Range("S9").Value = 50

# Row 10
$ws.Range("B10").Value = 56414
